# Generate Report for Handback
# Rename the "it-it" locale report to "de-de" and refresh the handoff/handback
# timestamps for that locale.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: rename the locale column header cell -----------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B1").Value = "de-de"

# Keep the backing table's column metadata in sync with the header cell.
$overviewTable = $overview.ListObjects.Item(1)
$overviewTable.ListColumns.Item(2).Name = "de-de"

# --- Locale sheet: rename "it-it" -> "de-de" -------------------------------
$localeSheet = $wb.Worksheets.Item("it-it")
$localeSheet.Name = "de-de"

# Rename the locale data table to match (Name also drives DisplayName).
$localeTable = $localeSheet.ListObjects.Item(1)
$localeTable.DisplayName = "de_de"
$localeTable.Name = "de-de"

# Refresh the handoff / handback datetimes for every row on the locale sheet.
$localeSheet.Range("E2:E5").Value = "2016-03-11 01:02:28"
$localeSheet.Range("H2:H5").Value = "2016-03-17 02:07:40"
